$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row is inserted at row 291 (pushing all subsequent rows
# down by one, up to the previous last row 426 which becomes 427).
$ws.Rows.Item(291).Insert()

# Seed the new row with the same row-level formatting/values as the row
# that is now directly below it (the old row 291, now shifted to row 292),
# then overwrite just the two cells that actually hold new data (the date
# and the volume).
$ws.Rows.Item(292).Copy()
$ws.Rows.Item(291).PasteSpecial()

$ws.Range("D291").Value = 44917
$ws.Range("J291").Value = 135
